# UniformA-HW20.xpc update ("New plotting functions, Matlab for Halfwidth
# data from Dr. Creuziger"):
#
#   - extend every existing row by one more HKL/pair column (W), value 21
#     in the numeric header row
#   - insert a new pair label "1Pair-B" into row 2's label list (so the
#     following labels shift one column to the right, and "MaxUnique"
#     now lands in the new W2 cell)
#   - append four new data rows (8-11) for the new half-width sampling
#     methods OffsetA, RD Single, TD Single, and HexGrid-90degTilt5degRes,
#     each filled with 1 across C:W

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W")

# ---------------------------------------------------------------------
# Row 1 : extend the numeric header from column V (20) to W (21)
# ---------------------------------------------------------------------
$ws.Range("V1").Copy() | Out-Null
$ws.Range("W1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("W1").Value = 21

# ---------------------------------------------------------------------
# Row 2 : insert "1Pair-B" into column O (right after "1Pair-A" in N),
#         shifting the remaining labels (old O2:V2 = "2Pairs-A" ...
#         "MaxUnique") one column to the right, into new P2:W2.
# ---------------------------------------------------------------------
$ws.Range("V2").Copy() | Out-Null
$ws.Range("W2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$shiftSource = @("O","P","Q","R","S","T","U","V")
$shiftTarget = @("P","Q","R","S","T","U","V","W")
$shiftValues = @()
foreach ($col in $shiftSource) {
    $shiftValues += , $ws.Range("${col}2").Value2
}
for ($i = $shiftValues.Length - 1; $i -ge 0; $i--) {
    $ws.Range("$($shiftTarget[$i])2").Value = $shiftValues[$i]
}
$ws.Range("O2").Value = "1Pair-B"

# ---------------------------------------------------------------------
# Rows 3-7 : existing sampling methods now also span the new column W
# ---------------------------------------------------------------------
foreach ($r in 3..7) {
    $ws.Range("V$r").Copy() | Out-Null
    $ws.Range("W$r").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range("W$r").Value = 1
}

# Row 7 becomes "OffsetF" - the 4 new methods are inserted ahead of the
# former last row ("HexGrid-90degTilt5degRes"), which itself moves down
# to become the new row 11.
$ws.Range("B7").Value = "OffsetF"

# ---------------------------------------------------------------------
# Rows 8-11 : four new sampling-method rows
# ---------------------------------------------------------------------
$newMethods = @("OffsetA", "RD Single", "TD Single", "HexGrid-90degTilt5degRes")

for ($i = 0; $i -lt $newMethods.Length; $i++) {
    $r = 8 + $i

    $ws.Range("A7").Copy() | Out-Null
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range("A$r").Value = 6 + $i

    $ws.Range("B$r").Value = $newMethods[$i]

    foreach ($col in $cols) {
        $ws.Range("$col$r").Value = 1
    }
}
